# CCC19 Derived Variables Spreadsheet
# Add new derived variable "O19a / VTE_comp_v2" as a new row of Table1,
# right after the existing "O19 / VTE_comp" row (row 95), pushing the
# remaining rows of the table down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new worksheet row at row 96 - this shifts every
# subsequent row (old 96..141) down by one (new 97..142) and keeps all
# of their existing values/formatting intact.
$ws.Rows.Item(96).Insert() | Out-Null

# Table1 currently covers A1:E141; grow it by one row so the table keeps
# wrapping the data (autoFilter / totals / styling all follow the range).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E142")) | Out-Null

# Populate the newly inserted row with the new variable, mirroring the
# columns used by its sibling row directly above it (Variable #, Variable
# Name (der_name), Category, Description).
$ws.Range("A96").Value = "O19a"
$ws.Range("B96").Value = "VTE_comp_v2"
$ws.Range("C96").Value = "Outcome"
$ws.Range("D96").Value = "Combined VTE complications (excluding SVT and thrombosis NOS)"

# Leave the view roughly where the author left it: scrolled near the new
# row, with that row selected.
$excel.ActiveWindow.ScrollRow = 84 | Out-Null
$ws.Range("A97").Select() | Out-Null
